$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Non-"Price" column text updates (Coin name, Link, Volume%) ---
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  -3.63%  '
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').Value = '  +1.46%  '
$ws.Range('E11').Value = '  -1.98%  '
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('E15').Value = '  -4.77%  '
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('E18').Value = '  +2.02%  '
$ws.Range('E19').Value = '  +1.47%  '
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('E21').Value = '  +1.61%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('E28').Value = '  +3.81%  '
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  +0.86%  '
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('E34').Value = '  +2.05%  '
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('E37').Value = '  +9.71%  '
$ws.Range('E38').Value = '  -3.62%  '
$ws.Range('E39').Value = '  +7.40%  '
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('E41').Value = '  +4.06%  '
$ws.Range('E42').Value = '  -4.12%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E44').Value = '  +6.04%  '
$ws.Range('E45').Value = '  +1.52%  '
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('E48').Value = '  +2.38%  '
$ws.Range('E49').Value = '  +4.34%  '
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('E51').Value = '  +0.28%  '

# --- "Price" column (D) updates: force text storage so values like
#     "25.506.85" or "0.04310" are not reinterpreted as numbers ---
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.506.85'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.664.54'
$ws.Range('D3').ClearFormats()
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9995'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '234.67'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4616'
$ws.Range('D7').ClearFormats()
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.662.99'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06933'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.59'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.337'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '74.85'
$ws.Range('D14').ClearFormats()
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5716'
$ws.Range('D15').ClearFormats()
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '25.513.69'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000006701'
$ws.Range('D19').ClearFormats()
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.36'
$ws.Range('D20').ClearFormats()
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.878.23'
$ws.Range('D21').ClearFormats()
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.644'
$ws.Range('D23').ClearFormats()
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.217'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '134.71'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '14.81'
$ws.Range('D26').ClearFormats()
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.709'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '103.56'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.947'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.07691'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.589'
$ws.Range('D32').ClearFormats()
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04310'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.620'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9399'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.5982'
$ws.Range('D36').ClearFormats()
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9170'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.477'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '105.75'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9997'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.828'
$ws.Range('D41').ClearFormats()
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.3701'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.973'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1111'
$ws.Range('D45').ClearFormats()
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.05254'
$ws.Range('D46').ClearFormats()
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.101'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '29.86'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.532'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.002'
$ws.Range('D50').ClearFormats()
